# Parameters_Abortion1_Preeclampsia05_EMM.xlsx edit
# "Explored different data generating parameters to get expected distributions."
#
# 1. SimParameters!B8 changes from 0.25 to 0.1 (the real data-generating-parameter
#    tweak described in the commit message). Every dependent formula elsewhere in
#    the workbook (potential_preec_trt!C19:F42, etc.) recalculates automatically.
# 2. postpreec_preg!C67:C83 and C108:C124 change from same-sheet back-references
#    (=C26 / =C27 filled down) to direct cross-sheet formulas
#    (=potential_preg_untrt!C<row>*SimParameters!$B$28), matching the pattern
#    already used in rows 26-42 / etc. The resulting values are unchanged.
# 3. View-state (selected cell) changes on a few sheets.

$wb = $excel.ActiveWorkbook

# --- 1. SimParameters!B8 0.25 -> 0.1 ------------------------------------------
$wsParams = $wb.Worksheets.Item("SimParameters")
$wsParams.Range("B8").Value = 0.1

# --- 2. postpreec_preg formulas: replace in-sheet references with direct refs -
$wsPost = $wb.Worksheets.Item("postpreec_preg")
for ($r = 67; $r -le 83; $r++) {
    $wsPost.Range("C$r").Formula = "=potential_preg_untrt!C$r*SimParameters!`$B`$28"
}
for ($r = 108; $r -le 124; $r++) {
    $wsPost.Range("C$r").Formula = "=potential_preg_untrt!C$r*SimParameters!`$B`$28"
}

# --- 3. Update the active-cell selections shown in the saved view state -------
# (Selecting a range on a sheet makes it the active sheet, so the sheets whose
# selection changed are visited first and SimParameters -- the sheet that must
# remain the active tab -- is re-selected last.)
$wb.Worksheets.Item("potential_preg_trt").Range("O37").Select()
$wb.Worksheets.Item("potential_preec_trt").Range("E24").Select()
$wb.Worksheets.Item("postpreec_preg").Range("I80").Select()
$wsParams.Range("C14").Select()

Write-Host "edit complete"
